$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2577906666666667
$ws.Range("H2").Value = 0.7733719999999999
$ws.Range("I2").Value = 0.1202607703685643
$ws.Range("J2").Value = 0.1202607703685642
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.04738633333333334
$ws.Range("N2").Value = 0.142159
$ws.Range("O2").Value = 0.05760194168856402
$ws.Range("P2").Value = 0.05760194168856402
$ws.Range("Q2").Value = 0.01221575446088889
$ws.Range("R2").Value = 0.109941790148
$ws.Range("S2").Value = 0.006927253882191826
$ws.Range("T2").Value = 0.006927253882191825
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2577906666666667
$ws.Range("H3").Value = 0.7733719999999999
$ws.Range("I3").Value = 0.1202607703685643
$ws.Range("J3").Value = 0.1202607703685642
$ws.Range("O3").Value = 0.7659981644722047
$ws.Range("P3").Value = 0.7659981644722047
$ws.Range("Q3").Value = 0.1624467026697778
$ws.Range("R3").Value = 1.462020324028
$ws.Range("S3").Value = 0.09211952936033352
$ws.Range("T3").Value = 0.0921195293603335
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2577906666666667
$ws.Range("H4").Value = 0.7733719999999999
$ws.Range("I4").Value = 0.1202607703685643
$ws.Range("J4").Value = 0.1202607703685642
$ws.Range("M4").Value = 0.1451156666666667
$ws.Range("N4").Value = 0.435347
$ws.Range("O4").Value = 0.1763998938392313
$ws.Range("P4").Value = 0.1763998938392313
$ws.Range("Q4").Value = 0.03740946445377778
$ws.Range("R4").Value = 0.336685180084
$ws.Range("S4").Value = 0.02121398712603891
$ws.Range("T4").Value = 0.02121398712603891
$ws.Range("I5").Value = 0.8797392296314358
$ws.Range("J5").Value = 0.8797392296314357
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04738633333333334
$ws.Range("N5").Value = 0.142159
$ws.Range("O5").Value = 0.05760194168856402
$ws.Range("P5").Value = 0.05760194168856402
$ws.Range("Q5").Value = 0.08936146330888889
$ws.Range("R5").Value = 0.8042531697800001
$ws.Range("S5").Value = 0.0506746878063722
$ws.Range("T5").Value = 0.05067468780637219
$ws.Range("I6").Value = 0.8797392296314358
$ws.Range("J6").Value = 0.8797392296314357
$ws.Range("O6").Value = 0.7659981644722047
$ws.Range("P6").Value = 0.7659981644722047
$ws.Range("S6").Value = 0.6738786351118712
$ws.Range("T6").Value = 0.6738786351118712
$ws.Range("I7").Value = 0.8797392296314358
$ws.Range("J7").Value = 0.8797392296314357
$ws.Range("M7").Value = 0.1451156666666667
$ws.Range("N7").Value = 0.435347
$ws.Range("O7").Value = 0.1763998938392313
$ws.Range("P7").Value = 0.1763998938392313
$ws.Range("Q7").Value = 0.2736600916377778
$ws.Range("S7").Value = 0.1551859067131924
$ws.Range("T7").Value = 0.1551859067131924
